$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "cocote" grading-scale column (I) used to cycle through A,B,C,D (4 options).
# This refactor drops option "D" from the scale, so the repeating pattern
# becomes A,B,C (3 options) -- every cell that used to read "D" now reads the
# next letter in the shortened cycle, and everything after it shifts up too.
$ws.Range("I6").Value  = "A"
$ws.Range("I7").Value  = "B"
$ws.Range("I8").Value  = "C"
$ws.Range("I9").Value  = "A"
$ws.Range("I10").Value = "B"
$ws.Range("I11").Value = "C"
$ws.Range("I12").Value = "A"
$ws.Range("I13").Value = "B"
$ws.Range("I14").Value = "C"
$ws.Range("I18").Value = "A"

# Move the active selection to B3 (was I22).
$ws.Range("B3").Select()
